# Add a "Salary" column to the Employees sheet, fix Mary's last name
# (Brown -> Jones), and populate the salary figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

# New header in D1
$ws.Range("D1").Value = "Salary"

# Correct Mary's last name on row 2
$ws.Range("B2").Value = "Jones"

# Salary values for each employee (rows 2-5)
$ws.Range("D2").Value = 200000
$ws.Range("D3").Value = 110000
$ws.Range("D4").Value = 135000
$ws.Range("D5").Value = 125000
